{"js": "// Implements the commit \"implemented view for player in Database\":\n//  1. Highlight the \"Spieler erstellen, bearbeiten, NICHT L\u00d6SCHEN!!!\" bullet\n//     in yellow (both the paragraph mark and the run).\n//  2. Move the \"_GoBack\" bookmark from the last bullet (\"Turnierverlauf\")\n//     to the end of the \"Spieler zum Verein hinzuf\u00fcgen, entfernen\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their (unique) text.\nlet highlightParagraph = null;\nlet bookmarkTargetParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"Spieler erstellen, bearbeiten, NICHT L\u00d6SCHEN!!!\") {\n    highlightParagraph = para;\n  } else if (text === \"Spieler zum Verein hinzuf\u00fcgen, entfernen\") {\n    bookmarkTargetParagraph = para;\n  }\n}\n\n// 1) Highlight the whole paragraph (run + paragraph mark) yellow.\nif (highlightParagraph) {\n  highlightParagraph.font.highlightColor = \"Yellow\";\n}\n\n// 2) Remove the existing \"_GoBack\" bookmark wherever it currently sits.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n\n// 3) Re-insert the \"_GoBack\" bookmark at the end of the target paragraph.\nif (bookmarkTargetParagraph) {\n  const endRange = bookmarkTargetParagraph.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Implements the commit \"implemented view for player in Database\":\n#  1. Highlight the \"Spieler erstellen, bearbeiten, NICHT LOESCHEN!!!\" bullet\n#     in yellow (both the run and the paragraph mark).\n#  2. Move the \"_GoBack\" bookmark from the last bullet (\"Turnierverlauf\")\n#     to the end of the \"Spieler zum Verein hinzufuegen, entfernen\" bullet.\n\n$d = $word.ActiveDocument\n\n$highlightParagraph = $null\n$bookmarkTargetParagraph = $null\n\nforeach ($p in $d.Paragraphs) {\n    # Paragraph.Range.Text includes the trailing paragraph-mark character.\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Spieler erstellen, bearbeiten, NICHT L\u00d6SCHEN!!!\") {\n        $highlightParagraph = $p\n    } elseif ($text -eq \"Spieler zum Verein hinzuf\u00fcgen, entfernen\") {\n        $bookmarkTargetParagraph = $p\n    }\n}\n\n# 1) Highlight the whole paragraph (run + paragraph mark) yellow.\nif ($highlightParagraph -ne $null) {\n    $highlightParagraph.Range.Font.HighlightColorIndex = 7   # wdYellow\n}\n\n# 2) Remove the existing \"_GoBack\" bookmark wherever it currently sits.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 3) Re-insert the \"_GoBack\" bookmark right after the text of the target\n#    paragraph (before its paragraph mark), matching the XML produced when\n#    Word places the bookmark immediately following the run.\nif ($bookmarkTargetParagraph -ne $null) {\n    $rng = $bookmarkTargetParagraph.Range.Duplicate\n    [void]$rng.MoveEnd(1, -1)    # wdCharacter: exclude the paragraph mark\n    $rng.Collapse(0)             # wdCollapseEnd: collapse to just after the text\n    $rng.InsertAfter([char]1)    # temporary placeholder so the range has length\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n    $rng.Text = \"\"               # remove the placeholder, bookmark stays collapsed here\n}\n"}
